$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2025-04-21 -> 2025-04-22, i.e. serial 45768 -> 45769) for every data
# row (rows 2 through 43).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45768) {
        $cell.Value2 = 45769
    }
}
